# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: Intel(R) Wi-Fi 6E AX211 160MHz driver version update + stat changes
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.160.0.4"
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 95.2

# Row 4: stat changes
$ws.Range("C4").Value = 627
$ws.Range("D4").Value = 97.90000000000001

# Row 5: stat change
$ws.Range("D5").Value = 98.40000000000001

# Row 6 (Totals row): Critical Minutes total
$ws.Range("C6").Value = 683

# Row 14: Total Samples update
$ws.Range("B14").Value = 449371

# Row 18: Total Samples update
$ws.Range("B18").Value = 77999
